# Applies the changes described by the commit:
#  - Sheet1 ("Hoja1"): C51 progress 50% -> 100%
#  - Sheet1: new row 70 -> task "Intro en ventana para crear cuota", responsible "Lucas"
#  - View state: scrolled down / selection moved to B71 (reflecting new data entry)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# C51: 0.5 -> 1 (task now complete)
$ws.Range("C51").Value = 1

# New row 70: task + responsible
$ws.Range("A70").Value = "Intro en ventana para crear cuota"
$ws.Range("B70").Value = "Lucas"

# Update the view to match where the user ended up editing
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 61
$ws.Range("B71").Select()
